$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: update date format (dd-mm-yyyy, unambiguous since day=28) and
# flip D3/G3 from 0 to 1
$ws.Range("A3").Value = "28-07-2022"
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

# Rows 4-21: only the date text changes (slashes -> dashes). Dates whose
# day-of-month is <= 12 are ambiguous to Excel's auto-detect (it may read
# them as a real date), so a leading apostrophe keeps them as literal text,
# matching the source data (plain inline strings).
$ws.Range("A4").Value = "'01-08-2022"
$ws.Range("A5").Value = "'04-08-2022"
$ws.Range("A6").Value = "'08-08-2022"
$ws.Range("A7").Value = "'11-08-2022"
$ws.Range("A8").Value = "15-08-2022"
$ws.Range("A9").Value = "18-08-2022"
$ws.Range("A10").Value = "22-08-2022"
$ws.Range("A11").Value = "25-08-2022"
$ws.Range("A12").Value = "29-08-2022"
$ws.Range("A13").Value = "'01-09-2022"
$ws.Range("A14").Value = "'05-09-2022"
$ws.Range("A15").Value = "'08-09-2022"
$ws.Range("A16").Value = "'12-09-2022"
$ws.Range("A17").Value = "15-09-2022"
$ws.Range("A18").Value = "19-09-2022"
$ws.Range("A19").Value = "22-09-2022"
$ws.Range("A20").Value = "26-09-2022"
$ws.Range("A21").Value = "29-09-2022"
